# Generate Report for Handoff
#
# The localization-status report is regenerated: the row for
# "dccb5d8a-...md" (already "In Translation") moves to row 2, and the
# row for "48404fab-...md" moves to row 3 and is updated from
# "In Translation" to "Ready for handoff" (new handoff xliff files,
# new handoff datetimes, Priority ht -> mt).

$wb = $excel.ActiveWorkbook

$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/843457f5855b9b45b839d4f862526f48e626ecfe/e2e/"
$url48404 = $urlBase + "48404fab-758b-4022-b734-91be03a04555.md"
$urlDccb5 = $urlBase + "dccb5d8a-6254-4a33-9903-cb30ecc4fdfb.md"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "dccb5d8a-6254-4a33-9903-cb30ecc4fdfb.md"
$ws.Range("A3").Value = "48404fab-758b-4022-b734-91be03a04555.md"

$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-22 14:13:42"

# Rebuild the two hyperlinks so row 2 points at rId2 (48404fab target,
# displaying the dccb5d8a path) and row 3 points at rId3 (dccb5d8a
# target, displaying the 48404fab path) - matching the swapped rows.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $url48404, [System.Type]::Missing, [System.Type]::Missing, "e2e\dccb5d8a-6254-4a33-9903-cb30ecc4fdfb.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $urlDccb5, [System.Type]::Missing, [System.Type]::Missing, "e2e\48404fab-758b-4022-b734-91be03a04555.md") | Out-Null

$ws.Columns.Item(5).ColumnWidth = 16.38265482584637
$ws.Columns.Item(6).ColumnWidth = 16.38265482584637

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "dccb5d8a-6254-4a33-9903-cb30ecc4fdfb.md"
$ws.Range("G2").Value = "dccb5d8a-6254-4a33-9903-cb30ecc4fdfb.4f9386eb9b277879a31d7fcdb842ab88f7ee0438.zh-cn.xlf"

$ws.Range("A3").Value = "48404fab-758b-4022-b734-91be03a04555.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("G3").Value = "48404fab-758b-4022-b734-91be03a04555.877ff5688decf5996d930e7c9c87891fb1950b04.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-22 14:13:34"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $url48404, [System.Type]::Missing, [System.Type]::Missing, "dccb5d8a-6254-4a33-9903-cb30ecc4fdfb.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $urlDccb5, [System.Type]::Missing, [System.Type]::Missing, "48404fab-758b-4022-b734-91be03a04555.md") | Out-Null

$ws.Columns.Item(3).ColumnWidth = 16.38265482584637

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "dccb5d8a-6254-4a33-9903-cb30ecc4fdfb.md"
$ws.Range("G2").Value = "dccb5d8a-6254-4a33-9903-cb30ecc4fdfb.4f9386eb9b277879a31d7fcdb842ab88f7ee0438.de-de.xlf"

$ws.Range("A3").Value = "48404fab-758b-4022-b734-91be03a04555.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("G3").Value = "48404fab-758b-4022-b734-91be03a04555.877ff5688decf5996d930e7c9c87891fb1950b04.de-de.xlf"
$ws.Range("H3").Value = "2016-08-22 14:13:42"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $url48404, [System.Type]::Missing, [System.Type]::Missing, "dccb5d8a-6254-4a33-9903-cb30ecc4fdfb.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $urlDccb5, [System.Type]::Missing, [System.Type]::Missing, "48404fab-758b-4022-b734-91be03a04555.md") | Out-Null

$ws.Columns.Item(3).ColumnWidth = 16.38265482584637
